# Defect table: add a new "Resolved on" column before the existing
# "Date of first occurence" column (old column D), shifting the old
# D:G columns to E:H. (Commit: "language of test report documents
# changed (to eng)" -> the sheet gained an English "Resolved on"
# header column.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at D; everything that was D:G (dates,
# time, describing test document, problem) slides right to E:H.
$ws.Columns("D").Insert()

# New header for the inserted column.
$ws.Range("D1").Value = "Resolved on"

# Give the new column a sensible width (close to the other header
# columns) - rows 2:5 are left blank under it, same as the source edit.
$ws.Columns("D").ColumnWidth = 11

# Match the author's final selection (cell D2, under the new header).
$ws.Range("D2").Select()
